# "Accept Button & Stuff"
# - Row 0: refreshed lichess game for row-0 challenge (new gameID/link/escrow rating).
# - Row 1: converted from a data row into the header row
#          (gameID, challenger, rating, wager, link, escrowID).
# - Rows 2-5: refreshed lichess games (new gameID/link/escrow rating), other
#             columns unchanged.
# - Rows 6-8: brand-new challenge rows appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 0: new game for this challenge row ----
$ws.Cells.Item(0,1).Value2 = "emGkGkKl"
$ws.Cells.Item(0,5).Value2 = "https://lichess.org/emGkGkKl"
$ws.Cells.Item(0,6).Value2 = 2317

# ---- Row 1: now the header row ----
$ws.Cells.Item(1,1).Value2 = "gameID"
$ws.Cells.Item(1,2).Value2 = "challenger"
$ws.Cells.Item(1,3).Value2 = "rating"
$ws.Cells.Item(1,4).Value2 = "wager"
$ws.Cells.Item(1,5).Value2 = "link"
$ws.Cells.Item(1,6).Value2 = "escrowID"

# ---- Row 2: new game ----
$ws.Cells.Item(2,1).Value2 = "p2LYrWLC"
$ws.Cells.Item(2,5).Value2 = "https://lichess.org/p2LYrWLC"
$ws.Cells.Item(2,6).Value2 = 2559

# ---- Row 3: new game ----
$ws.Cells.Item(3,1).Value2 = "oIOLTnRN"
$ws.Cells.Item(3,5).Value2 = "https://lichess.org/oIOLTnRN"
$ws.Cells.Item(3,6).Value2 = 2560

# ---- Row 4: new game ----
$ws.Cells.Item(4,1).Value2 = "NHiX7o3n"
$ws.Cells.Item(4,5).Value2 = "https://lichess.org/NHiX7o3n"
$ws.Cells.Item(4,6).Value2 = 2561

# ---- Row 5: new game ----
$ws.Cells.Item(5,1).Value2 = "YVTmakze"
$ws.Cells.Item(5,5).Value2 = "https://lichess.org/YVTmakze"
$ws.Cells.Item(5,6).Value2 = 2562

# ---- Row 6: brand-new challenge row ----
$ws.Cells.Item(6,1).Value2 = "gcDcW4K8"
$ws.Cells.Item(6,2).Value2 = "trashboatsr"
$ws.Cells.Item(6,3).Value2 = 1818
$ws.Cells.Item(6,4).Value2 = 100
$ws.Cells.Item(6,5).Value2 = "https://lichess.org/gcDcW4K8"
$ws.Cells.Item(6,6).Value2 = 2564

# ---- Row 7: brand-new challenge row ----
$ws.Cells.Item(7,1).Value2 = "8dQmNTlC"
$ws.Cells.Item(7,2).Value2 = "trashboatsr"
$ws.Cells.Item(7,3).Value2 = 1818
$ws.Cells.Item(7,4).Value2 = 100
$ws.Cells.Item(7,5).Value2 = "https://lichess.org/8dQmNTlC"
$ws.Cells.Item(7,6).Value2 = 2565

# ---- Row 8: brand-new challenge row ----
$ws.Cells.Item(8,1).Value2 = "532PWje3"
$ws.Cells.Item(8,2).Value2 = "trashboatsr"
$ws.Cells.Item(8,3).Value2 = 1818
$ws.Cells.Item(8,4).Value2 = 100
$ws.Cells.Item(8,5).Value2 = "https://lichess.org/532PWje3"
$ws.Cells.Item(8,6).Value2 = 2566
